$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 (cohort 2023, period_index 3): num_customers 38 -> 39, retention_rate recalculated
$ws.Range("C31").Value = 39
$ws.Range("E31").Value = 0.01686851211072665

# Row 34 (cohort 2024, period_index 2): num_customers 58 -> 61, retention_rate recalculated
$ws.Range("C34").Value = 61
$ws.Range("E34").Value = 0.02703900709219858

# Row 37 (cohort 2025, period_index 0): num_customers 588 -> 599, cohort_size 588 -> 599
$ws.Range("C37").Value = 599
$ws.Range("D37").Value = 599
